# Spain Segunda update (19-06-2024 21:51)
# The source feed re-sorted same-kickoff-time fixtures, which produced a
# swap of the data (everything except the running "id" index in column A)
# between several pairs of adjacent rows, plus a few odds corrections on
# the very last row (no swap partner available there).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $rowA, $rowB, $firstCol, $lastCol) {
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)
        $valA = $cellA.Value2
        $valB = $cellB.Value2
        $cellA.Value2 = $valB
        $cellB.Value2 = $valA
    }
}

# Columns B (2) through AD (30) get swapped; column A (the sequential id)
# stays put on both rows.
$firstCol = 2
$lastCol = 30

$rowPairs = @(
    @(83, 84),
    @(96, 97),
    @(103, 104),
    @(142, 143),
    @(169, 170),
    @(190, 191),
    @(213, 214),
    @(228, 229),
    @(235, 236),
    @(253, 254),
    @(289, 290),
    @(300, 301),
    @(446, 447)
)

foreach ($pair in $rowPairs) {
    Swap-Rows $ws $pair[0] $pair[1] $firstCol $lastCol
}

# Row 459 is the last row in the sheet (no partner row to swap with); its
# odds were simply corrected in place.
$ws.Range("O459").Value2 = 1.909
$ws.Range("Q459").Value2 = 3.8
$ws.Range("S459").Value2 = 1.975
$ws.Range("T459").Value2 = 1.875
$ws.Range("V459").Value2 = 2.1
$ws.Range("W459").Value2 = 1.775
